$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 983.7692
$ws.Range("I17").Value = 700
$ws.Range("J17").Value = 1007.4167
$ws.Range("K17").Value = 2100
$ws.Range("L17").Value = 3022.2501
$ws.Range("M17").Value = -1932
$ws.Range("N17").Value = -3358.2501
$ws.Range("H43").Value = 769
$ws.Range("I43").Value = 525.25
$ws.Range("J43").Value = 826.35297
$ws.Range("K43").Value = 525.25
$ws.Range("L43").Value = 826.35297
$ws.Range("M43").Value = -456.25
$ws.Range("N43").Value = -964.35297
$ws.Range("H58").Value = 1409.3846
$ws.Range("I58").Value = 402
$ws.Range("J58").Value = 2584.6667
$ws.Range("K58").Value = 1206
$ws.Range("L58").Value = 7754.000100000001
$ws.Range("M58").Value = -1056
$ws.Range("N58").Value = -8054.000100000001
$ws.Range("H62").Value = 1288.8334
$ws.Range("I62").Value = 1081.6666
$ws.Range("J62").Value = 1496
$ws.Range("K62").Value = 1081.6666
$ws.Range("L62").Value = 1496
$ws.Range("M62").Value = -457.6666
$ws.Range("N62").Value = -2744
$ws.Range("H65").Value = 1288.8334
$ws.Range("I65").Value = 1081.6666
$ws.Range("J65").Value = 1496
$ws.Range("K65").Value = 5408.333000000001
$ws.Range("L65").Value = 7480
$ws.Range("M65").Value = -2288.333000000001
$ws.Range("N65").Value = -13720
$ws.Range("H76").Value = 3021.9167
$ws.Range("I76").Value = 2826.3
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 2826.3
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2511.3
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 3021.9167
$ws.Range("I79").Value = 2826.3
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 2826.3
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -1734.3
$ws.Range("N79").Value = -6184
$ws.Range("H132").Value = 2234.0386
$ws.Range("I132").Value = 1122.8235
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 3368.4705
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -838.4704999999999
$ws.Range("N132").Value = -18059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1686.9565
$ws.Range("I102").Value = 1575
$ws.Range("J102").Value = 2433.3333
$ws.Range("K102").Value = 1575
$ws.Range("L102").Value = 2433.3333
$ws.Range("M102").Value = 47
$ws.Range("N102").Value = -5677.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1455.5555
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 1366.6666
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 1366.6666
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -3612.6666
$ws.Range("H89").Value = 1455.5555
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 1366.6666
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 6833.333000000001
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -18065.333
$ws.Range("H105").Value = 1811.2963
$ws.Range("I105").Value = 1558.4706
$ws.Range("J105").Value = 2241.1
$ws.Range("K105").Value = 1558.4706
$ws.Range("L105").Value = 2241.1
$ws.Range("M105").Value = 188.5293999999999
$ws.Range("N105").Value = -5735.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1513.7368
$ws.Range("I58").Value = 872.6
$ws.Range("J58").Value = 1742.7142
$ws.Range("K58").Value = 872.6
$ws.Range("L58").Value = 1742.7142
$ws.Range("M58").Value = -669.6
$ws.Range("N58").Value = -2148.7142
$ws.Range("H62").Value = 68610.664
$ws.Range("I62").Value = 201262
$ws.Range("J62").Value = 2285
$ws.Range("K62").Value = 201262
$ws.Range("L62").Value = 2285
$ws.Range("M62").Value = -200638
$ws.Range("N62").Value = -3533
$ws.Range("H65").Value = 68610.664
$ws.Range("I65").Value = 201262
$ws.Range("J65").Value = 2285
$ws.Range("K65").Value = 1006310
$ws.Range("L65").Value = 11425
$ws.Range("M65").Value = -1003190
$ws.Range("N65").Value = -17665
$ws.Range("H105").Value = 2614.2856
$ws.Range("I105").Value = 500
$ws.Range("J105").Value = 4200
$ws.Range("K105").Value = 500
$ws.Range("L105").Value = 4200
$ws.Range("M105").Value = 1247
$ws.Range("H132").Value = 1790.1111
$ws.Range("I132").Value = 1207.6316
$ws.Range("J132").Value = 3173.5
$ws.Range("K132").Value = 3622.8948
$ws.Range("L132").Value = 9520.5
$ws.Range("M132").Value = -1092.8948
$ws.Range("N132").Value = -14580.5
$ws.Range("H136").Value = 1513.7368
$ws.Range("I136").Value = 872.6
$ws.Range("J136").Value = 1742.7142
$ws.Range("K136").Value = 2617.8
$ws.Range("L136").Value = 5228.142599999999
$ws.Range("M136").Value = -67.80000000000018
$ws.Range("N136").Value = -10328.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 905.55554
$ws.Range("I131").Value = 287.66666
$ws.Range("J131").Value = 982.7917
$ws.Range("K131").Value = 862.9999799999999
$ws.Range("L131").Value = 2948.3751
$ws.Range("M131").Value = 4177.00002
$ws.Range("N131").Value = -13028.3751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5914.839
$ws.Range("I70").Value = 4938.4
$ws.Range("J70").Value = 9983.333000000001
$ws.Range("K70").Value = 4938.4
$ws.Range("L70").Value = 9983.333000000001
$ws.Range("M70").Value = -4668.4
$ws.Range("N70").Value = -10523.333
$ws.Range("H73").Value = 5914.839
$ws.Range("I73").Value = 4938.4
$ws.Range("J73").Value = 9983.333000000001
$ws.Range("K73").Value = 4938.4
$ws.Range("L73").Value = 9983.333000000001
$ws.Range("M73").Value = -4002.4
$ws.Range("N73").Value = -11855.333
$ws.Range("H75").Value = 42875
$ws.Range("I75").Value = 44714.285
$ws.Range("J75").Value = 30000
$ws.Range("K75").Value = 44714.285
$ws.Range("L75").Value = 30000
$ws.Range("M75").Value = -43840.285
$ws.Range("N75").Value = -31748
$ws.Range("H78").Value = 42875
$ws.Range("I78").Value = 44714.285
$ws.Range("J78").Value = 30000
$ws.Range("K78").Value = 134142.855
$ws.Range("L78").Value = 90000
$ws.Range("M78").Value = -129774.855
$ws.Range("N78").Value = -98736
$ws.Range("H97").Value = 2240
$ws.Range("I97").Value = 1487.25
$ws.Range("J97").Value = 2842.2
$ws.Range("K97").Value = 1487.25
$ws.Range("L97").Value = 2842.2
$ws.Range("M97").Value = -991.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3350.25
$ws.Range("I100").Value = 2934.3333
$ws.Range("J100").Value = 3488.889
$ws.Range("K100").Value = 2934.3333
$ws.Range("L100").Value = 3488.889
$ws.Range("M100").Value = -2393.3333
$ws.Range("N100").Value = -4570.889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3312.875
$ws.Range("I62").Value = 3083.8333
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 3083.8333
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -2459.8333
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 3312.875
$ws.Range("I65").Value = 3083.8333
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 15419.1665
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -12299.1665
$ws.Range("N65").Value = -26240
$ws.Range("H96").Value = 2542.8572
$ws.Range("I96").Value = 1600
$ws.Range("J96").Value = 3250
$ws.Range("K96").Value = 1600
$ws.Range("L96").Value = 3250
$ws.Range("M96").Value = -227
$ws.Range("N96").Value = -5996
